# Insert a new row for "Table 2" at row 2, shifting existing data down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2 (pushes existing row 2+ down by one).
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the Table 2 entry.
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = "No code"

# Update the active selection to A3, matching the post-edit workbook state.
$ws.Range("A3").Select()
